$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) updates: force text format first so numeric-looking
# strings (e.g. "128.33", "1.00", "7.30") are not silently coerced to
# Double values (which would round/truncate and change the cell type),
# then restore the default "Normal" style so no stray formatting is left
# behind on cells that started out with the default style.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '54.289.78'
$ws.Range("E2").Value = '  +0.80%  '

Set-TextValue "D3" '2.278.09'
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("E4").Value = '  +0.34%  '

Set-TextValue "D5" '497.89'
$ws.Range("E5").Value = '  +1.48%  '

Set-TextValue "D6" '128.33'
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("E7").Value = '  +0.43%  '

Set-TextValue "D8" '0.528'
$ws.Range("E8").Value = '  +0.10%  '

Set-TextValue "D9" '0.0954'
$ws.Range("E9").Value = '  +2.40%  '

$ws.Range("E10").Value = '  +1.25%  '

$ws.Range("E11").Value = '  +3.17%  '

$ws.Range("E12").Value = '  +1.08%  '

Set-TextValue "D13" '2.683.67'
$ws.Range("E13").Value = '  +0.56%  '

Set-TextValue "D14" '22.59'
$ws.Range("E14").Value = '  +5.32%  '

Set-TextValue "D15" '54.239.68'
$ws.Range("E15").Value = '  +0.79%  '

$ws.Range("E16").Value = '  +0.48%  '

Set-TextValue "D17" '2.280.75'
$ws.Range("E17").Value = '  +1.54%  '

Set-TextValue "D18" '10.23'
$ws.Range("E18").Value = '  +4.52%  '

Set-TextValue "D19" '4.12'
$ws.Range("E19").Value = '  +2.16%  '

Set-TextValue "D20" '304.44'
$ws.Range("E20").Value = '  +2.33%  '

$ws.Range("E21").Value = '  +2.78%  '

Set-TextValue "D22" '1.00'
$ws.Range("E22").Value = '  +0.42%  '

Set-TextValue "D23" '61.82'
$ws.Range("E23").Value = '  -2.88%  '

Set-TextValue "D24" '0.998'
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("E25").Value = '  +2.11%  '

Set-TextValue "D26" '7.30'
$ws.Range("E26").Value = '  +2.61%  '

$ws.Range("E27").Value = '  +8.13%  '

$ws.Range("E28").Value = '  +0.96%  '

Set-TextValue "D29" '5.95'
$ws.Range("E29").Value = '  +2.44%  '

$ws.Range("E30").Value = '  +0.75%  '

$ws.Range("E31").Value = '  +1.46%  '

$ws.Range("E32").Value = '  +0.20%  '

Set-TextValue "D33" '17.75'
$ws.Range("E33").Value = '  +1.82%  '

Set-TextValue "D34" '0.997'
$ws.Range("E34").Value = '  +0.05%  '

Set-TextValue "D35" '0.924'
$ws.Range("E35").Value = '  +10.33%  '

$ws.Range("E36").Value = '  +1.12%  '

Set-TextValue "D37" '3.74'
$ws.Range("E37").Value = '  +2.89%  '

$ws.Range("E38").Value = '  -0.69%  '

$ws.Range("E39").Value = '  +1.33%  '

$ws.Range("E40").Value = '  +1.50%  '

Set-TextValue "D41" '125.48'
$ws.Range("E41").Value = '  -0.31%  '

Set-TextValue "D42" '4.76'
$ws.Range("E42").Value = '  -2.84%  '

$ws.Range("E43").Value = '  +2.66%  '

Set-TextValue "D44" '0.0897'
$ws.Range("E44").Value = '  +0.78%  '

Set-TextValue "D46" '239.77'
$ws.Range("E46").Value = '  -1.07%  '

$ws.Range("E47").Value = '  -0.56%  '

$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("E49").Value = '  +1.14%  '

Set-TextValue "D50" '16.27'
$ws.Range("E50").Value = '  +0.56%  '

$ws.Range("E51").Value = '  +0.36%  '
